$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Amount" values for the AHUs (I/J columns) and Pumps (L/M columns)
$ws.Range("J2").Value = 2
$ws.Range("M2").Value = 2

$ws.Range("J3").Value = 1
$ws.Range("M3").Value = 1

$ws.Range("J4").Value = 2
$ws.Range("M4").Value = 2

# Update the active selection to J5
$ws.Range("J5").Select()
